# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and for two pairs of rows whose ranking order changed, also update
# Coin (B) and Link (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = $null; C = $null; D = '42.123.19'; E = '  -3.66%  ' }
    @{ Row = 3; B = $null; C = $null; D = '2.193.72'; E = '  -3.82%  ' }
    @{ Row = 4; B = $null; C = $null; D = $null; E = '  +1.08%  ' }
    @{ Row = 5; B = $null; C = $null; D = '105.21'; E = '  -15.14%  ' }
    @{ Row = 6; B = $null; C = $null; D = '290.80'; E = '  +9.19%  ' }
    @{ Row = 7; B = $null; C = $null; D = '0.615'; E = '  -3.24%  ' }
    @{ Row = 8; B = $null; C = $null; D = '1.01'; E = '  +0.04%  ' }
    @{ Row = 9; B = $null; C = $null; D = '0.585'; E = '  -6.54%  ' }
    @{ Row = 10; B = $null; C = $null; D = '43.28'; E = '  -10.67%  ' }
    @{ Row = 11; B = $null; C = $null; D = '0.0896'; E = '  -5.49%  ' }
    @{ Row = 12; B = $null; C = $null; D = '53.97'; E = '  -0.67%  ' }
    @{ Row = 13; B = $null; C = $null; D = '8.58'; E = '  -7.25%  ' }
    @{ Row = 14; B = $null; C = $null; D = '0.103'; E = '  -3.22%  ' }
    @{ Row = 15; B = $null; C = $null; D = '0.911'; E = '  +0.79%  ' }
    @{ Row = 16; B = $null; C = $null; D = '14.65'; E = '  -5.07%  ' }
    @{ Row = 17; B = $null; C = $null; D = '2.529.09'; E = '  -3.66%  ' }
    @{ Row = 18; B = $null; C = $null; D = '2.247.21'; E = '  -1.37%  ' }
    @{ Row = 19; B = $null; C = $null; D = '42.261.41'; E = '  -3.39%  ' }
    @{ Row = 20; B = $null; C = $null; D = '7.09'; E = '  +1.11%  ' }
    @{ Row = 21; B = $null; C = $null; D = '0.0000103'; E = '  -6.14%  ' }
    @{ Row = 22; B = $null; C = $null; D = '72.30'; E = '  -0.15%  ' }
    @{ Row = 23; B = $null; C = $null; D = '3.34'; E = '  +15.84%  ' }
    @{ Row = 24; B = $null; C = $null; D = '2.23'; E = '  -8.86%  ' }
    @{ Row = 25; B = $null; C = $null; D = '223.99'; E = '  -4.93%  ' }
    @{ Row = 26; B = $null; C = $null; D = '8.95'; E = '  -5.32%  ' }
    @{ Row = 27; B = $null; C = $null; D = '0.997'; E = '  -2.02%  ' }
    @{ Row = 28; B = $null; C = $null; D = '11.38'; E = '  -3.84%  ' }
    @{ Row = 29; B = $null; C = $null; D = '3.88'; E = '  -0.63%  ' }
    @{ Row = 30; B = $null; C = $null; D = $null; E = '  -2.15%  ' }
    @{ Row = 31; B = 'WEMIXToken'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '3.20'; E = '  -5.01%  ' }
    @{ Row = 32; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '37.33'; E = '  -12.81%  ' }
    @{ Row = 33; B = $null; C = $null; D = '172.17'; E = '  -0.32%  ' }
    @{ Row = 34; B = $null; C = $null; D = '20.52'; E = '  -5.53%  ' }
    @{ Row = 35; B = $null; C = $null; D = '0.0856'; E = '  -6.43%  ' }
    @{ Row = 36; B = $null; C = $null; D = '5.45'; E = '  -5.48%  ' }
    @{ Row = 37; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '4.73'; E = '  +1.50%  ' }
    @{ Row = 38; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '4.13'; E = '  -3.47%  ' }
    @{ Row = 39; B = $null; C = $null; D = '0.124'; E = '  -4.78%  ' }
    @{ Row = 40; B = $null; C = $null; D = '0.0355'; E = '  -5.82%  ' }
    @{ Row = 41; B = $null; C = $null; D = '0.100'; E = '  -6.40%  ' }
    @{ Row = 42; B = $null; C = $null; D = '2.39'; E = '  -5.90%  ' }
    @{ Row = 43; B = $null; C = $null; D = '0.227'; E = '  -5.60%  ' }
    @{ Row = 44; B = $null; C = $null; D = '68.68'; E = '  -8.49%  ' }
    @{ Row = 45; B = $null; C = $null; D = '1.01'; E = '  +0.54%  ' }
    @{ Row = 46; B = $null; C = $null; D = '12.30'; E = '  -12.15%  ' }
    @{ Row = 47; B = $null; C = $null; D = '1.27'; E = '  -7.55%  ' }
    @{ Row = 48; B = $null; C = $null; D = '5.32'; E = '  -5.42%  ' }
    @{ Row = 49; B = $null; C = $null; D = '1.26'; E = '  -0.37%  ' }
    @{ Row = 50; B = $null; C = $null; D = '100.69'; E = '  -1.26%  ' }
    @{ Row = 51; B = $null; C = $null; D = '8.35'; E = '  -3.24%  ' }
)

foreach ($item in $data) {
    $r = $item.Row

    if ($null -ne $item.B) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($null -ne $item.C) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($null -ne $item.D) {
        # Price column is stored as literal text in the workbook (values like
        # "42.123.19" or "1.01"). Force the cell to Text format before writing
        # so Excel doesn't auto-convert the string into a number, then drop
        # back to the Normal style so no stray formatting is left behind.
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($null -ne $item.E) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}

Write-Output "Updated $($data.Count) rows in cryptos sheet"
